$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.780.56"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.144.39"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "575.25"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "148.41"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.141.50"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "36.99"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "3.660.82"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "64.913.11"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "3.148.44"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "500.46"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").Value = "15.14"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D25").Value = "83.62"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "8.80"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  +6.18%  "
$ws.Range("D31").Value = "27.40"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "6.13"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "54.53"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "0.0887"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("D38").Value = "476.52"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").Value = "2.95"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "3.004.12"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  -3.70%  "
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "27.98"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D50").Value = "2.22"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "33.20"
$ws.Range("E51").Value = "  +7.19%  "
